# Add "Lore of Light" spells to the Spells sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spells")

$rows = @(
  @("Cursed Blades","New",1,"Enchantment","No","No","No","No","No","No","No","Yes","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Decrepify","New",4,"Necromancy","No","No","Yes","No","No","No","Yes","Yes","Yes","1.0.0","Complete","Publicly Released","Not on website"),
  @("Guiding Wind","New",1,"Enchantment","No","No","Yes","Yes","No","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Incantation of Protection","New",2,"Abjuration","No","Yes","Yes","No","Yes","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Incantation of Zeal","New",2,"Abjuration","No","Yes","Yes","No","Yes","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Skull Storm","New",6,"Necromancy","No","No","Yes","No","No","No","No","Yes","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Tomb Strike","New",0,"Necromancy","No","No","Yes","No","No","No","No","Yes","Yes","1.0.0","Complete","Publicly Released","Not on website"),
  @("Burning Gaze","New",1,"Evocation","No","No","Yes","No","No","No","No","No","No","2.0.0","Complete","Publicly Released","Not on website"),
  @("Dazzling Light","New",3,"Evocation","No","Yes","Yes","No","No","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Holy Nova","New",4,"Evocation","No","No","Yes","No","Yes","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Mass Haste","New",8,"Trasmutation","No","No","No","No","No","No","Yes","No","Yes","1.0.0","Complete","Publicly Released","Not on website"),
  @("Net of Light","New",4,"Conjuration","No","No","Yes","No","Yes","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
  @("Warding Light","New",1,"Abjuration","No","No","Yes","No","Yes","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website")
)

$startRow = 63
for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $startRow + $i
  $row = $rows[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

# Match the author's final view/selection state
$ws.Activate()
$ws.Range("D67").Select()

$wb.Save()
